# Update "想去人数" (interest count) figures in the "展览" and "全部类型"
# sheets to reflect newly generated output.

$wb = $excel.ActiveWorkbook

# Map of cell address -> new value, shared by both sheets that list the
# same exhibition rows (展览 holds exhibition-only rows; 全部类型 aggregates
# every category, so the exhibition rows there land on slightly different
# row numbers further down the sheet).
$updatesSheet1 = @{
    "F2"  = 622
    "F3"  = 497
    "F4"  = 1296
    "F5"  = 1169
    "F6"  = 14337
    "F7"  = 16531
    "F9"  = 96
    "F12" = 201
    "F14" = 50
    "F19" = 104
    "F20" = 37
    "F21" = 1264
    "F24" = 38
    "F27" = 6720
    "F28" = 971
    "F30" = 1119
    "F31" = 11
    "F33" = 5748
    "F34" = 103
    "F36" = 190
    "F37" = 4820
    "F38" = 18
}

$updatesSheet4 = @{
    "F2"  = 622
    "F3"  = 497
    "F4"  = 1296
    "F5"  = 1169
    "F6"  = 14337
    "F7"  = 16531
    "F9"  = 96
    "F12" = 201
    "F14" = 50
    "F19" = 104
    "F20" = 37
    "F21" = 1264
    "F25" = 38
    "F28" = 6720
    "F29" = 971
    "F31" = 1119
    "F32" = 11
    "F36" = 5748
    "F37" = 103
    "F39" = 190
    "F40" = 4820
    "F41" = 18
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($addr in $updatesSheet1.Keys) {
    $ws1.Range($addr).Value = $updatesSheet1[$addr]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($addr in $updatesSheet4.Keys) {
    $ws4.Range($addr).Value = $updatesSheet4[$addr]
}
